$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 117; this pushes the existing rows
# 117-189 down to 118-190 (dimension grows from R189 to R190).
$ws.Rows("117:117").Insert()

# Populate the newly inserted row 117 with the new price-report entry.
$ws.Range("A117").Value = 10
$ws.Range("B117").Value = "Vega Modelo de Temuco"
$ws.Range("C117").Value = "La Araucanía"
$ws.Range("D117").Value = 44322
$ws.Range("E117").Value = 9
$ws.Range("F117").Value = 100112052
$ws.Range("G117").Value = "Albahaca"
$ws.Range("H117").Value = "Sin especificar"
$ws.Range("I117").Value = "Primera"
$ws.Range("J117").Value = 15
$ws.Range("K117").Value = 8000
$ws.Range("L117").Value = 8000
$ws.Range("M117").Value = 8000
$ws.Range("N117").Value = "`$/docena de matas"
$ws.Range("O117").Value = "Región de Arica y Parinacota"
$ws.Range("P117").Value = 1333
$ws.Range("Q117").Value = 6
$ws.Range("R117").Value = "Hortaliza"
